$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '67.574.72'
$ws.Range('E2').Value = '  +4.16%  '

$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '3.256.15'
$ws.Range('E3').Value = '  +3.47%  '

$ws.Range('E4').Value = '  +0.04%  '

$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '578.53'
$ws.Range('E5').Value = '  +2.06%  '

$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '181.03'
$ws.Range('E6').Value = '  +6.18%  '

$ws.Range('E7').Value = '  -0.04%  '

$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.592'
$ws.Range('E8').Value = '  -3.98%  '

$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '3.254.28'
$ws.Range('E9').Value = '  +3.60%  '

$ws.Range('E10').Value = '  +5.22%  '

$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '6.79'
$ws.Range('E11').Value = '  +3.72%  '

$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.415'
$ws.Range('E12').Value = '  +5.67%  '

$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '3.824.71'
$ws.Range('E13').Value = '  +3.66%  '

$ws.Range('E14').Value = '  +1.59%  '

$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '28.65'
$ws.Range('E15').Value = '  +6.27%  '

$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '67.530.74'
$ws.Range('E16').Value = '  +4.26%  '

$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '0.0000168'
$ws.Range('E17').Value = '  +3.12%  '

$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '3.257.34'
$ws.Range('E18').Value = '  +3.72%  '

$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '5.83'
$ws.Range('E19').Value = '  +2.71%  '

$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '13.56'
$ws.Range('E20').Value = '  +5.93%  '

$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '376.92'
$ws.Range('E21').Value = '  +6.28%  '

$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '7.64'
$ws.Range('E22').Value = '  +5.46%  '

$ws.Range('E23').Value = '  -0.08%  '

$ws.Range('E24').Value = '  +4.21%  '

$ws.Range('E25').Value = '  +2.70%  '

$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '0.0000120'
$ws.Range('E26').Value = '  +2.46%  '

$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '9.64'
$ws.Range('E27').Value = '  +1.15%  '

$ws.Range('E28').Value = '  +3.52%  '

$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '1.01'
$ws.Range('E29').Value = '  +0.50%  '

$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '5.77'
$ws.Range('E30').Value = '  +7.87%  '

$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '1.98'
$ws.Range('E31').Value = '  +4.20%  '

$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '22.64'
$ws.Range('E32').Value = '  +3.40%  '

$ws.Range('E33').Value = '  -0.01%  '

$ws.Range('E34').Value = '  +6.64%  '

$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '6.93'
$ws.Range('E35').Value = '  +4.40%  '

$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '163.65'
$ws.Range('E36').Value = '  +6.49%  '

$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '1.50'
$ws.Range('E37').Value = '  +4.39%  '

$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.853'
$ws.Range('E38').Value = '  +2.79%  '

$ws.Range('E39').Value = '  +6.46%  '

$ws.Range('B40').Value = 'EnergySwap'
$ws.Range('C40').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '26.85'
$ws.Range('E40').Value = '  +3.00%  '

$ws.Range('B41').Value = 'RenderToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '6.79'
$ws.Range('E41').Value = '  +12.90%  '

$ws.Range('B42').Value = 'Filecoin'
$ws.Range('C42').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '4.56'
$ws.Range('E42').Value = '  +9.49%  '

$ws.Range('B43').Value = 'dogwifhat'
$ws.Range('C43').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '2.62'
$ws.Range('E43').Value = '  +4.96%  '

$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '365.42'
$ws.Range('E44').Value = '  +12.72%  '

$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '2.742.65'
$ws.Range('E45').Value = '  +3.10%  '

$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '25.43'
$ws.Range('E46').Value = '  +4.88%  '

$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '40.90'
$ws.Range('E47').Value = '  +4.19%  '

$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '0.0678'
$ws.Range('E48').Value = '  +3.21%  '

$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '0.0280'
$ws.Range('E49').Value = '  +2.68%  '

$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '1.00'

$ws.Range('E51').Value = '  +0.10%  '
